# Applies the "Tried constrained PRIM, works" edit:
#   - CHP1 "time" row: Min Value 4692.553052 -> 5375.493591
#   - CHP2 block: "heat_pump" and "time" rows swap places (time now first),
#     with the "time" Min Value updated 4539.922962 -> 5227.175586
#   - CHP3 block: reordered from (heat_pump, duration_increase, time) to
#     (time, heat_pump, duration_increase), with the "time" Min Value
#     updated 4540.717407 -> 5373.471119
#   - Selection moves from A5 to D9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay a text cell even though it looks
# like a number (matches how the source file stores these as shared
# strings rather than numeric cells).
function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
}

# --- CHP1: just a value tweak on the existing "time" row ---
Set-TextValue "C4" "5375.493591"

# --- CHP2: swap "heat_pump" (row 6) and "time" (row 7) rows ---
# Row 6 becomes the "time" row with an updated Min Value
$ws.Range("B6").Value = "time"
Set-TextValue "C6" "5227.175586"
Set-TextValue "D6" "5997.172835"
$ws.Range("E6").ClearContents()

# Row 7 becomes the "heat_pump" row
$ws.Range("B7").Value = "heat_pump"
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E7").Value = "{True}       {True}"

# --- CHP3: reorder heat_pump / duration_increase / time -> time / heat_pump / duration_increase ---
# Row 8 becomes the "time" row with an updated Min Value
$ws.Range("B8").Value = "time"
Set-TextValue "C8" "5373.471119"
Set-TextValue "D8" "5999.593822"
$ws.Range("E8").ClearContents()

# Row 9 becomes the "heat_pump" row
$ws.Range("B9").Value = "heat_pump"
$ws.Range("C9").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("E9").Value = "{True}       {True}"

# Row 10 becomes the "duration_increase" row
$ws.Range("B10").Value = "duration_increase"
Set-TextValue "C10" "1000"
Set-TextValue "D10" "1000"

# --- Selection moves to D9 ---
$ws.Range("D9").Select()
